$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
    'schwingen',
    'posten',
    'töten',
    'heulen',
    'geben',
    'kümmern',
    'kosten',
    'wollen',
    'kürzen',
    'spüren',
    'tropfen',
    'fassen',
    'fahren',
    'warnen',
    'pflanzen',
    'heben',
    'bleiben',
    'flehen',
    'knien',
    'kichern',
    'zeigen',
    'klagen',
    'filmen',
    'klingen',
    'trennen',
    'kratzen',
    'wundern',
    'nerven',
    'bergen',
    'münzen',
    'bitten',
    'sagen',
    'suchen',
    'bluten',
    'wachsen',
    'siegen',
    'mögen',
    'leeren',
    'enden',
    'wärmen',
    'malen',
    'dringen',
    'wirken',
    'jagen',
    'fesseln',
    'lockern',
    'dienen',
    'machen',
    'heilen',
    'jubeln',
    'fließen',
    'kriegen',
    'leiden',
    'seufzen',
    'wüten',
    'spinnen',
    'fangen',
    'reizen',
    'nähen',
    'bieten',
    'schlucken',
    'meinen',
    'mauern',
    'sinken',
    'schenken',
    'spielen',
    'helfen',
    'rasen',
    'stoppen',
    'bauen',
    'drohen',
    'zeichnen',
    'hacken',
    'schreiten',
    'krachen',
    'quälen',
    'schützen',
    'formen',
    'flüchten',
    'achten',
    'boxen',
    'schmecken',
    'dürfen',
    'brauchen',
    'fühlen',
    'hauen',
    'knarren',
    'äußern',
    'stärken',
    'pfeifen',
    'stecken',
    'ärgern',
    'planen',
    'platzen',
    'küssen',
    'läuten',
    'spenden',
    'fischen',
    'führen',
    'träumen',
    'lügen',
    'sichern',
    'gelten',
    'bellen',
    'folgen',
    'bremsen',
    'prügeln',
    'wehtun',
    'klingeln',
    'gründen',
    'schwimmen',
    'tollen',
    'schwächen',
    'scheinen',
    'scheitern',
    'sperren',
    'strahlen',
    'kennen',
    'spannen',
    'schrecken',
    'biegen',
    'irren',
    'stürmen',
    'altern',
    'freuen',
    'rufen',
    'klettern',
    'arten',
    'wagen',
    'zünden',
    'stehlen',
    'passen',
    'tragen',
    'zielen',
    'grüßen',
    'pissen',
    'sterben',
    'kehren',
    'stimmen',
    'runden',
    'backen',
    'graben',
    'wechseln',
    'erben',
    'feiern',
    'zögern',
    'sprengen',
    'ändern',
    'ehren',
    'narren',
    'trauen',
    'weichen',
    'leisten',
    'saufen',
    'stammen',
    'greifen',
    'treiben',
    'buchen',
    'sorgen',
    'schulden',
    'rücken',
    'drehen',
    'sammeln',
    'füttern',
    'danken',
    'wenden',
    'liegen',
    'werfen',
    'lesen',
    'wandern',
    'treten',
    'starten',
    'streichen',
    'fallen',
    'löschen',
    'tauchen',
    'schalten',
    'fällen',
    'hören',
    'liefern',
    'regnen',
    'tauschen',
    'werden',
    'füllen',
    'hindern',
    'pressen',
    'schwören',
    'räumen',
    'betteln',
    'decken',
    'teilen',
    'loben'
)

for ($i = 0; $i -lt $words.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $words[$i]
}

